# Replace SVN with Git
#
# The source document was authored with SVN-era wording ("... in
# subversion and the resulting build of the master branch ...").  This
# script updates that sentence to reference Git / the "main" branch
# instead, matching the rest of the repository's modern tooling.

$d = $word.ActiveDocument

# The document was saved with Track Changes turned on (w:trackRevisions
# in settings.xml). Make sure our edits land as plain text changes
# rather than tracked insertions/deletions.
$d.TrackRevisions = $false

# "... a tag is created from this commit in subversion and the
# resulting build ..."  ->  "... in Git and the resulting build ..."
$d.Content.Find.Execute(
    "subversion", $true, $false, $false, $false, $false,
    $true, 1, $false, "Git", 2) | Out-Null

# "... of the master branch ..."  ->  "... of the main branch ..."
$d.Content.Find.Execute(
    "master branch", $true, $false, $false, $false, $false,
    $true, 1, $false, "main branch", 2) | Out-Null
